# School_Details workbook update
# - Row 22 (MGA School Kenya / MGA Kenya / kenyamoe1) had its Raw_SectionMap
#   text sitting in the Class_Sections column (F22) with nothing in
#   Raw_SectionMap (G22). This fills F22 with the proper human-readable
#   Class_Sections text and moves the existing Raw_SectionMap-style text
#   into G22, matching every other row's layout.
# - Widens column F so the new multi-line text is readable.
# - Leaves the view scrolled/selected near the edited cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The value that used to (incorrectly) live in F22 is the Raw_SectionMap
# style text - preserve it so it can be moved into G22.
$rawSectionMap = $ws.Range("F22").Value

$classSections = "06: 06A0 (06A0)(Grade 6)`r`n07: 07A0 (07A0)(Grade 7)`r`n08: 08A0 (08A0)(Grade 8)`r`n09: 09A0 (09A0)(Grade 9)`r`n10: 10A0 (10A0)(Grade 10)`r`n11: 11A0 (11A0)(Grade 10)`r`n12: 12A0 (12A0)(Grade 10)`r`nGrade 1: 01A0 (Grade 1)`r`nGrade 2: 02A0 (Grade 2)`r`nGrade 3: 03A0 (Grade 3)`r`nGrade 4: 04A0 (Grade 4)`r`nGrade 5: 05A0 (Grade 5)"

# Put the new Class_Sections text in F22, and the original text into the
# (previously empty) G22, mirroring the pattern used by every other row.
$ws.Range("F22").Value = $classSections
$ws.Range("G22").Value = $rawSectionMap

# Setting a multi-line value auto-sizes the row; re-fit the row so it goes
# back to the workbook's implicit (default) height, same as all other rows.
$ws.Range("F22:G22").EntireRow.AutoFit()

# Widen column F to comfortably fit the new multi-line Class_Sections text.
$ws.Columns.Item(6).ColumnWidth = 44.85

# Scroll the view down and leave the selection near the edited area.
$ws.Range("A17").Select()
$ws.Range("F28").Select()
